$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E16").Value = "2307"
$ws.Range("E17").Value = "2306"
$ws.Range("E18").Value = "2305"
$ws.Range("E19").Value = "2304"
$ws.Range("E20").Value = "2303"
$ws.Range("E21").Value = "2302"
$ws.Range("E22").Value = "2301"

$ws.Range("F16").Value = 34666
$ws.Range("F22").Value = 40000
